$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the existing date/time cell formats down onto the new rows so the
#    new cells reuse the same style indices (s="1" for dates, s="2" for
#    times) instead of minting brand-new cellXfs entries.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range($ws.Cells.Item(8, 1), $ws.Cells.Item(28, 1)).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Range($ws.Cells.Item(8, 2), $ws.Cells.Item(19, 2)).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Dates / times for the new rows (A = date serial, B = time fraction).
#    Rows 20-28 only ever had a date (no time-of-day value).
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = 44105
$ws.Cells.Item(8, 2).Value = 0.4548611111111111
$ws.Cells.Item(9, 1).Value = 44105
$ws.Cells.Item(9, 2).Value = 0.4548611111111111
$ws.Cells.Item(10, 1).Value = 44105
$ws.Cells.Item(10, 2).Value = 0.51041666666666663
$ws.Cells.Item(11, 1).Value = 44105
$ws.Cells.Item(11, 2).Value = 0.51041666666666663
$ws.Cells.Item(12, 1).Value = 44105
$ws.Cells.Item(12, 2).Value = 0.55902777777777779
$ws.Cells.Item(13, 1).Value = 44105
$ws.Cells.Item(13, 2).Value = 0.55902777777777779

$ws.Cells.Item(14, 1).Value = 44208
$ws.Cells.Item(14, 2).Value = 0.35416666666666669
$ws.Cells.Item(15, 1).Value = 44208
$ws.Cells.Item(15, 2).Value = 0.40625
$ws.Cells.Item(16, 1).Value = 44208
$ws.Cells.Item(16, 2).Value = 0.46875

$ws.Cells.Item(17, 1).Value = 44214
$ws.Cells.Item(17, 2).Value = 0.37847222222222227
$ws.Cells.Item(18, 1).Value = 44214
$ws.Cells.Item(18, 2).Value = 0.41319444444444442
$ws.Cells.Item(19, 1).Value = 44214
$ws.Cells.Item(19, 2).Value = 0.44444444444444442

$ws.Cells.Item(20, 1).Value = 44319
$ws.Cells.Item(21, 1).Value = 44319
$ws.Cells.Item(22, 1).Value = 44319
$ws.Cells.Item(23, 1).Value = 44319
$ws.Cells.Item(24, 1).Value = 44319
$ws.Cells.Item(25, 1).Value = 44319
$ws.Cells.Item(26, 1).Value = 44319
$ws.Cells.Item(27, 1).Value = 44319
$ws.Cells.Item(28, 1).Value = 44319

# ---------------------------------------------------------------------------
# 3) Site column (C) for the new rows 8-28, written top-to-bottom.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 3).Value = "IC-Z1"
$ws.Cells.Item(9, 3).Value = "IC-Z2"
$ws.Cells.Item(10, 3).Value = "IC-C1"
$ws.Cells.Item(11, 3).Value = "IC-C2"
$ws.Cells.Item(12, 3).Value = "IC-U1"
$ws.Cells.Item(13, 3).Value = "IC-U2"
$ws.Cells.Item(14, 3).Value = "IC-U1"
$ws.Cells.Item(15, 3).Value = "IC-U2"
$ws.Cells.Item(16, 3).Value = "IC-C1"
$ws.Cells.Item(17, 3).Value = "IC-C2"
$ws.Cells.Item(18, 3).Value = "IC-Z1"
$ws.Cells.Item(19, 3).Value = "IC-Z2"
$ws.Cells.Item(20, 3).Value = "IC-Z2"
$ws.Cells.Item(21, 3).Value = "IC-Z3"
$ws.Cells.Item(22, 3).Value = "IC-Z1"
$ws.Cells.Item(23, 3).Value = "IC-C2"
$ws.Cells.Item(24, 3).Value = "IC-C3"
$ws.Cells.Item(25, 3).Value = "IC-C1"
$ws.Cells.Item(26, 3).Value = "IC-U2"
$ws.Cells.Item(27, 3).Value = "IC-U3"
$ws.Cells.Item(28, 3).Value = "IC-U1"

# ---------------------------------------------------------------------------
# 4) Series column (D). The previous "100-112" .. "600-612" placeholders are
#    replaced with the new "101-112-1" .. "601-612-1" labels, then every new
#    row gets its series label too.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = "101-112-1"
$ws.Cells.Item(3, 4).Value = "201-212-1"
$ws.Cells.Item(4, 4).Value = "301-312-1"
$ws.Cells.Item(5, 4).Value = "401-412-1"
$ws.Cells.Item(6, 4).Value = "501-512-1"
$ws.Cells.Item(7, 4).Value = "601-612-1"

$ws.Cells.Item(8, 4).Value = "501-512-2"
$ws.Cells.Item(9, 4).Value = "601-612-2"
$ws.Cells.Item(10, 4).Value = "301-312-2"
$ws.Cells.Item(11, 4).Value = "401-412-2"
$ws.Cells.Item(12, 4).Value = "101-112-2"
$ws.Cells.Item(13, 4).Value = "201-212-2"

$ws.Cells.Item(14, 4).Value = "101-112-3"
$ws.Cells.Item(15, 4).Value = "201-212-3"
$ws.Cells.Item(16, 4).Value = "301-312-3"

$ws.Cells.Item(17, 4).Value = "401-412-3"
$ws.Cells.Item(18, 4).Value = "501-512-3"
$ws.Cells.Item(19, 4).Value = "601-612-3"

$ws.Cells.Item(20, 4).Value = "601-612-4"
$ws.Cells.Item(22, 4).Value = "501-512-4"
$ws.Cells.Item(21, 4).Value = "901-912-1"
$ws.Cells.Item(23, 4).Value = "401-412-4"
$ws.Cells.Item(24, 4).Value = "801-812-1"
$ws.Cells.Item(25, 4).Value = "301-312-4"
$ws.Cells.Item(26, 4).Value = "201-212-4"
$ws.Cells.Item(27, 4).Value = "701-712-1"
$ws.Cells.Item(28, 4).Value = "101-112-4"

# ---------------------------------------------------------------------------
# 5) Match the author's final selection/active cell.
# ---------------------------------------------------------------------------
$ws.Range("C28").Select() | Out-Null
